# Update supplier/brand names in the discount table (rows 11-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: supplier "Gb corp" -> "GB Corp", brand "lassa" -> "Lassa"
$ws.Range("A11").Value = "GB Corp"
$ws.Range("C11").Value = "Lassa"

# Row 12: supplier "Gb corp" -> "GB Corp", brand "Yokahama" -> "Yokohama"
$ws.Range("A12").Value = "GB Corp"
$ws.Range("C12").Value = "Yokohama"

# Row 13: supplier "Gb corp" -> "GB Corp", brand "SunFull" -> "Sunfull"
$ws.Range("A13").Value = "GB Corp"
$ws.Range("C13").Value = "Sunfull"

# Row 14: supplier "Gb corp" -> "GB Corp" (brand "Thunderer" unchanged)
$ws.Range("A14").Value = "GB Corp"

# Row 15: supplier "Gb corp" -> "GB Corp", brand "DoubleCoin" -> "Double Coin"
$ws.Range("A15").Value = "GB Corp"
$ws.Range("C15").Value = "Double Coin"

# Row 16: supplier "gresco" -> "Gresco" (brand "Kumho" unchanged)
$ws.Range("A16").Value = "Gresco"

# Row 17: supplier "gresco" -> "Gresco", brand "Wanda" -> "Winda"
$ws.Range("A17").Value = "Gresco"
$ws.Range("C17").Value = "Winda"

# The retyped brand cells (column C, rows 11-17) lost their larger font
# formatting and fell back to the default style/font size.
$ws.Range("C11:C17").Font.Size = 11

# Update the active selection to reflect where the user ended up working
$ws.Range("A17").Select()
